# Projet de visualition de donnee One Piece Explorer
# Convert the "Ventes au Japon" column from text (formatted with thousands
# separators, stored as shared strings) to real numeric values, and fix the
# "Annee" header (drop the accent).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Header row
$ws.Range("A1").Value = "Annee"
$ws.Range("B1").Value = "Ventes au Japon"

# Yearly sales figures, previously text like "1,822,218", now plain numbers.
# Index 0 -> row 2 (year 1997) ... index 24 -> row 26 (year 2021).
$years = @(1997,1998,1999,2000,2001,2002,2003,2004,2005,2006,2007,2008,2009,2010,2011,2012,2013,2014,2015,2016,2017,2018,2019,2020,2021)
$sales = @(1822218,5310427,5580271,5355206,4813183,4254073,3794620,3581786,3382467,3434777,3444014,4261054,5002885,5307870,4734778,3782159,3799410,4089586,4442492,3747273,2752727,2868164,2988003,2630293,2357214)

for ($i = 0; $i -lt $years.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $years[$i]
    $ws.Cells.Item($row, 2).Value = $sales[$i]
}

$ws.Range("A1").Select() | Out-Null
